$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) to reflect the new "through" date
$ws.Name = "Through 2022-02-23"

# Update the header label in I1 (shared string) to match the new date
$ws.Range("I1").Value = "2022 (through 02-23)"

# Update the February (row 3) total for the "through" column
$ws.Range("I3").Value = 118

# Update the grand Total row (row 14) for the "through" column
$ws.Range("I14").Value = 277

$wb.Save()
